# Updated cryptos list values (Price and Volume(1h) columns) per diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.407.07"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.309.22"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.11%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.49"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.70"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.01%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.309.08"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.467"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.84"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.54%  "

$ws.Range("E11").Value = "  -3.44%  "

$ws.Range("E12").Value = "  -1.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.876.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.06%  "

$ws.Range("E14").Value = "  +0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.86"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.310.91"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.75%  "

$ws.Range("E17").Value = "  -2.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.408.03"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.12"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.18"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.55%  "

$ws.Range("E21").Value = "  -2.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.35%  "

$ws.Range("E23").Value = "  -0.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.537"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.446.94"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.06%  "

$ws.Range("E27").Value = "  -7.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.171"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.18"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.88%  "

$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("E32").Value = "  -3.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.58"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.64"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.68%  "

$ws.Range("E35").Value = "  -6.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.13"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.57"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.51%  "

$ws.Range("E38").Value = "  -3.32%  "

$ws.Range("E39").Value = "  -2.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.338.71"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.75"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -13.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0729"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.95"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.67%  "

$ws.Range("E44").Value = "  -3.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.15"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.76%  "

$ws.Range("E46").Value = "  -5.03%  "

$ws.Range("E47").Value = "  -4.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.363.23"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -7.00%  "

$ws.Range("E49").Value = "  +0.01%  "

$ws.Range("E50").Value = "  -6.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.31"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.77%  "
